# ProjectPlan.docx edit script
# Implements the changes described in the commit "Updated QM, added risklist + in process worklist":
#  1. "Android Studios" -> "Android Studio" in the software bullet list (first occurrence only).
#  2. "These risk are of varying types" -> "These risks are of varying types".
#  3. Italicize "team members not being committed to the project".
#  4. Replace "React Native" with "Android Studio" (both mentions inside the technical-risk
#     paragraph) - the team actually used Arduino + Android Studio, not React Native.
#  5. Refresh the stale cached page-number field in the footer (2 -> 8) to reflect the
#     document's growth after the edits above.

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2 (Word constants used positionally below)

# --- 1. "Android Studios" -> "Android Studio" (only the first mention, in the bullet list) ---
$rng = $d.Content
$rng.Find.Execute("Android Studios", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Android Studio", 1) | Out-Null

# --- 2. Fix the grammar slip "These risk are" -> "These risks are" ---
$rng = $d.Content
$rng.Find.Execute("These risk are of varying types", $false, $false, $false, $false, $false, $true, 1, $false, `
    "These risks are of varying types", 2) | Out-Null

# --- 3. Italicize "team members not being committed to the project" ---
$rng = $d.Content
$found = $rng.Find.Execute("team members not being committed to the project", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Italic = 1
}

# --- 4. "React Native" -> "Android Studio" (both quoted mentions) ---
$rng = $d.Content
$rng.Find.Execute("React Native", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Android Studio", 2) | Out-Null

# --- 5. Update the stale cached PAGE field result in the footer (2 -> 8) ---
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    foreach ($f in $ftr.Range.Fields) {
        $result = $f.Result
        if ($result.Text -eq "2") {
            $result.Characters.Item(1).Text = "8"
        }
    }
}
